$wb = $excel.ActiveWorkbook

# ----- Sheet "展览" (sheet1) -----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 111
$ws1.Range("F5").Value = 1723
$ws1.Range("F6").Value = 3297
$ws1.Range("F7").Value = 968
$ws1.Range("F8").Value = 2141
$ws1.Range("F9").Value = 2061
$ws1.Range("F10").Value = 1072
$ws1.Range("F11").Value = 572
$ws1.Range("F13").Value = 1643
$ws1.Range("F14").Value = 362
$ws1.Range("F15").Value = 74
$ws1.Range("F16").Value = 26
$ws1.Range("F17").Value = 82
$ws1.Range("F18").Value = 148
$ws1.Range("F19").Value = 1513
$ws1.Range("F20").Value = 572
$ws1.Range("F21").Value = 673
$ws1.Range("F22").Value = 560
$ws1.Range("F23").Value = 12002
$ws1.Range("F24").Value = 12019
$ws1.Range("F25").Value = 884
$ws1.Range("F26").Value = 677
$ws1.Range("F28").Value = 5
$ws1.Range("F29").Value = 293
$ws1.Range("F30").Value = 1882
$ws1.Range("F31").Value = 174
$ws1.Range("F32").Value = 510

# ----- Sheet "本地生活" (sheet3) -----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 69

# ----- Sheet "全部类型" (sheet4) -----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 69
$ws4.Range("F6").Value = 111
$ws4.Range("F7").Value = 1723
$ws4.Range("F8").Value = 3297
$ws4.Range("F9").Value = 968
$ws4.Range("F10").Value = 2141
$ws4.Range("F11").Value = 2061
$ws4.Range("F12").Value = 1072
$ws4.Range("F13").Value = 572
$ws4.Range("F15").Value = 1643
$ws4.Range("F16").Value = 362
$ws4.Range("F17").Value = 74
$ws4.Range("F18").Value = 26
$ws4.Range("F20").Value = 82
$ws4.Range("F22").Value = 148
$ws4.Range("F23").Value = 1513
$ws4.Range("F24").Value = 572
$ws4.Range("F25").Value = 673
$ws4.Range("F26").Value = 560
$ws4.Range("F27").Value = 12002
$ws4.Range("F28").Value = 12019
$ws4.Range("F29").Value = 884
$ws4.Range("F30").Value = 677
$ws4.Range("F32").Value = 5
$ws4.Range("F33").Value = 293
$ws4.Range("F34").Value = 1882
$ws4.Range("F37").Value = 174
$ws4.Range("F38").Value = 510
